{"js": "// Update the methods/abstract text: survey years and sampling-location counts.\nconst body = context.document.body;\n\n// 1) \"...between March and August at 24 sampling locations selected at random (12 in Low...\"\n//    -> \"...between March and August in 2001 and 2002 at 29 sampling locations selected at random (13 in Low...\"\nconst hit1 = body.search(\n  \"between March and August at 24 sampling locations selected at random (12 in Low\",\n  { matchCase: true }\n);\nhit1.load(\"items\");\nawait context.sync();\nif (hit1.items.length > 0) {\n  hit1.items[0].insertText(\n    \"between March and August in 2001 and 2002 at 29 sampling locations selected at random (13 in Low\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 2) \"Low Flow Channel, 12 in High\" -> \"Low Flow Channel, 16 in High\"\nconst hit2 = body.search(\"Low Flow Channel, 12 in High\", { matchCase: true });\nhit2.load(\"items\");\nawait context.sync();\nif (hit2.items.length > 0) {\n  hit2.items[0].insertText(\"Low Flow Channel, 16 in High\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Re-type \"network of StowAway electronic thermistors.\" (same text) so the stray\n//    spell-check proofErr markers around \"StowAway\" collapse away, matching a clean\n//    re-save. Insert the replacement text immediately before the match, then delete\n//    the original range (instead of a plain Replace) so the proofErr siblings that sit\n//    between the original runs are removed along with the rest of the old range.\nconst hit3 = body.search(\"network of StowAway electronic thermistors.\", { matchCase: true });\nhit3.load(\"items\");\nawait context.sync();\nif (hit3.items.length > 0) {\n  const target = hit3.items[0];\n  target.insertText(\"network of StowAway electronic thermistors.\", Word.InsertLocation.before);\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# Update the methods/abstract text: survey years and sampling-location counts.\n$d = $word.ActiveDocument\n\n# 1) \"...between March and August at 24 sampling locations selected at random (12 in Low...\"\n#    -> \"...between March and August in 2001 and 2002 at 29 sampling locations selected at random (13 in Low...\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"between March and August at 24 sampling locations selected at random (12 in Low\"\n$find.Replacement.Text = \"between March and August in 2001 and 2002 at 29 sampling locations selected at random (13 in Low\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) \"Low Flow Channel, 12 in High\" -> \"Low Flow Channel, 16 in High\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Low Flow Channel, 12 in High\"\n$find2.Replacement.Text = \"Low Flow Channel, 16 in High\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n\n# 3) Re-type \"StowAway electronic thermistors.\" (same text) so the stray spell-check\n#    proofErr markers around \"StowAway\" collapse into a single run, matching a clean re-save.\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = \"network of StowAway electronic thermistors.\"\n$find3.Replacement.Text = \"network of StowAway electronic thermistors.\"\n$find3.Execute($find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2) | Out-Null\n"}
